$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.428.58"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.652.41"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.17%  "
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.27"
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = "  +0.69%  "
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.01"
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.20%  "
$origStyle_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("D8").Style = $origStyle_D8
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "2.663.43"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  +9.50%  "
$ws.Range("E11").Value = "  -2.55%  "
$origStyle_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("D12").Style = $origStyle_D12
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D14").Value = "3.116.72"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "59.425.21"
$ws.Range("E15").Value = "  -0.05%  "
$origStyle_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.10"
$ws.Range("D16").Style = $origStyle_D16
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "2.649.34"
$ws.Range("E18").Value = "  -0.05%  "
$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.65"
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = "  -2.93%  "
$origStyle_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.43"
$ws.Range("D20").Style = $origStyle_D20
$ws.Range("E20").Value = "  -1.86%  "
$origStyle_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.34"
$ws.Range("D21").Style = $origStyle_D21
$ws.Range("E21").Value = "  -0.24%  "
$origStyle_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.29"
$ws.Range("D22").Style = $origStyle_D22
$ws.Range("E22").Value = "  +1.16%  "
$origStyle_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.995"
$ws.Range("D23").Style = $origStyle_D23
$ws.Range("E23").Value = "  -0.57%  "
$origStyle_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.63"
$ws.Range("D24").Style = $origStyle_D24
$ws.Range("E24").Value = "  +2.54%  "
$ws.Range("E25").Value = "  +1.49%  "
$origStyle_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.412"
$ws.Range("D26").Style = $origStyle_D26
$ws.Range("E26").Value = "  -1.26%  "
$origStyle_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = $origStyle_D27
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "0.0₃0803"
$ws.Range("E28").Value = "  -0.76%  "
$origStyle_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("D29").Style = $origStyle_D29
$ws.Range("E29").Value = "  -0.22%  "
$origStyle_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.68"
$ws.Range("D30").Style = $origStyle_D30
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +0.46%  "
$origStyle_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.76"
$ws.Range("D33").Style = $origStyle_D33
$ws.Range("E33").Value = "  -1.04%  "
$origStyle_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.50"
$ws.Range("D34").Style = $origStyle_D34
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  +2.07%  "
$origStyle_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.20"
$ws.Range("D36").Style = $origStyle_D36
$ws.Range("E36").Value = "  +2.24%  "
$origStyle_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.897"
$ws.Range("D37").Style = $origStyle_D37
$ws.Range("E37").Value = "  -4.61%  "
$origStyle_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.883"
$ws.Range("D38").Style = $origStyle_D38
$ws.Range("E38").Value = "  +1.30%  "
$origStyle_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.91"
$ws.Range("D39").Style = $origStyle_D39
$ws.Range("E39").Value = "  +0.59%  "
$origStyle_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("D40").Style = $origStyle_D40
$ws.Range("E40").Value = "  +1.69%  "
$origStyle_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.59"
$ws.Range("D41").Style = $origStyle_D41
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("E42").Value = "  +4.17%  "
$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.46"
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$origStyle_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.77"
$ws.Range("D45").Style = $origStyle_D45
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E46").Value = "  -1.86%  "
$origStyle_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0537"
$ws.Range("D47").Style = $origStyle_D47
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "2.051.80"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("E49").Value = "  +2.02%  "
$origStyle_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.79"
$ws.Range("D50").Style = $origStyle_D50
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$origStyle_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.10"
$ws.Range("D51").Style = $origStyle_D51
$ws.Range("E51").Value = "  +1.09%  "
